# Apply corrections and rewordings described by the commit diff.
$d = $word.ActiveDocument

# --- Spelling fixes ---
$d.Content.Find.Execute("Reacions", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Reactions", 2) | Out-Null

$d.Content.Find.Execute("messgaes", $true, $false, $false, $false, $false,
                         $true, 1, $false, "messages", 2) | Out-Null

$d.Content.Find.Execute("assosiated", $true, $false, $false, $false, $false,
                         $true, 1, $false, "associated", 2) | Out-Null

# --- Content rewordings ---
$d.Content.Find.Execute(
    "Migrating one-on-one conversations or direct messages from Slack to Teams, ensuring that private communication is transferred securely.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Direct messages from Slack cannot be fully migrated to Teams. Some private conversations may not be transferred due to technical limitations.",
    2) | Out-Null

$d.Content.Find.Execute(
    "Transferring user groups or teams from Slack to Teams, preserving the group structure and membership for seamless collaboration.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "User groups from Slack cannot be fully transferred to Teams. Group structure and membership may not be preserved during migration.",
    2) | Out-Null

$d.Content.Find.Execute(
    "Transferring regular messages from Slack to Teams, including text-based communication which user sent to himself.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Self messages (messages sent by users to themselves) cannot be migrated from Slack to Teams.",
    2) | Out-Null
